# edit.ps1 -- PowerShell-style PowerPoint COM-interop script
#
# Reproduces the two functional changes described by the commit diff:
#
#   1. On slide 5, the table's table-style reference is switched from the
#      deck's custom style ({ED5542EF-3B5B-4BF1-A826-76CDAA587BB2}, still
#      defined in ppt/tableStyles.xml) to the built-in style
#      {00573AD5-7961-479D-856E-4788A1573543}.
#
#   2. The presentation's theme (ppt/theme/theme1.xml, the theme used by
#      the slide master / slides) is switched from the "Integral" / "Red
#      Violet" palette to the standard "Office" palette (the palette that,
#      before this edit, only the Notes Master's theme part used). The
#      font scheme and format scheme are already identical between the
#      two themes in this deck, so only the 12 theme colors actually need
#      to change.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: convert an "RRGGBB" hex string into the BGR-packed long that
# PowerPoint's ColorFormat.RGB / VBA RGB() expects.
# ---------------------------------------------------------------------
function HexToRgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536 + $g * 256 + $r)
}

# ---------------------------------------------------------------------
# 1) Table style on slide 5's table.
# ---------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{00573AD5-7961-479D-856E-4788A1573543}")
    }
}

# ---------------------------------------------------------------------
# 2) Recolor the main theme (theme1.xml) from "Red Violet" to "Office".
#    Index order of ThemeColorScheme: dk1, lt1, dk2, lt2, accent1-6,
#    hlink, folHlink.
# ---------------------------------------------------------------------
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$mainTheme = $p.SlideMaster.Theme
$mainColorScheme = $mainTheme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $mainColorScheme.Item($i).RGB = (HexToRgbLong $officeColors[$i - 1])
}
